$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 130
$ws.Cells.Item(4, 9).Value = 36.2
$ws.Cells.Item(4, 11).Value = 36.2
$ws.Cells.Item(4, 13).Value = 77.8
$ws.Cells.Item(28, 8).Value = 874.9167
$ws.Cells.Item(28, 9).Value = 435.2
$ws.Cells.Item(28, 10).Value = 1607.7778
$ws.Cells.Item(28, 11).Value = 435.2
$ws.Cells.Item(28, 12).Value = 1607.7778
$ws.Cells.Item(28, 13).Value = 49.80000000000001
$ws.Cells.Item(28, 14).Value = -2577.7778
$ws.Cells.Item(62, 8).Value = 9748.5
$ws.Cells.Item(62, 10).Value = 13323
$ws.Cells.Item(62, 12).Value = 13323
$ws.Cells.Item(62, 14).Value = -14571
$ws.Cells.Item(65, 8).Value = 9748.5
$ws.Cells.Item(65, 10).Value = 13323
$ws.Cells.Item(65, 12).Value = 66615
$ws.Cells.Item(65, 14).Value = -72855
$ws.Cells.Item(98, 8).Value = 2402.3333
$ws.Cells.Item(98, 9).Value = 2324.3845
$ws.Cells.Item(98, 10).Value = 2529
$ws.Cells.Item(98, 11).Value = 2324.3845
$ws.Cells.Item(98, 12).Value = 2529
$ws.Cells.Item(98, 13).Value = -826.3845000000001
$ws.Cells.Item(98, 14).Value = -5525
$ws.Cells.Item(122, 8).Value = 2402.3333
$ws.Cells.Item(122, 9).Value = 2324.3845
$ws.Cells.Item(122, 10).Value = 2529
$ws.Cells.Item(122, 11).Value = 6973.1535
$ws.Cells.Item(122, 12).Value = 7587
$ws.Cells.Item(122, 13).Value = -4523.1535
$ws.Cells.Item(122, 14).Value = -12487
$ws.Cells.Item(137, 8).Value = 14394.879
$ws.Cells.Item(137, 9).Value = 5641.143
$ws.Cells.Item(137, 11).Value = 16923.429
$ws.Cells.Item(137, 13).Value = -14373.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1928.7441
$ws.Cells.Item(2, 9).Value = 650.30304
$ws.Cells.Item(2, 10).Value = 6147.6
$ws.Cells.Item(2, 11).Value = 650.30304
$ws.Cells.Item(2, 12).Value = 6147.6
$ws.Cells.Item(2, 13).Value = -537.30304
$ws.Cells.Item(2, 14).Value = -6373.6
$ws.Cells.Item(45, 8).Value = 10885.333
$ws.Cells.Item(45, 9).Value = 11262.4
$ws.Cells.Item(45, 11).Value = 11262.4
$ws.Cells.Item(45, 13).Value = -10885.4
$ws.Cells.Item(74, 8).Value = 15435.258
$ws.Cells.Item(74, 9).Value = 2562.0588
$ws.Cells.Item(74, 10).Value = 31067
$ws.Cells.Item(74, 11).Value = 2562.0588
$ws.Cells.Item(74, 12).Value = 31067
$ws.Cells.Item(74, 13).Value = -1688.0588
$ws.Cells.Item(74, 14).Value = -32815
$ws.Cells.Item(77, 8).Value = 15435.258
$ws.Cells.Item(77, 9).Value = 2562.0588
$ws.Cells.Item(77, 10).Value = 31067
$ws.Cells.Item(77, 11).Value = 12810.294
$ws.Cells.Item(77, 12).Value = 155335
$ws.Cells.Item(77, 13).Value = -8442.293999999998
$ws.Cells.Item(77, 14).Value = -164071
$ws.Cells.Item(116, 8).Value = 1928.7441
$ws.Cells.Item(116, 9).Value = 650.30304
$ws.Cells.Item(116, 10).Value = 6147.6
$ws.Cells.Item(116, 11).Value = 650.30304
$ws.Cells.Item(116, 12).Value = 6147.6
$ws.Cells.Item(116, 13).Value = 1643.69696
$ws.Cells.Item(116, 14).Value = -10735.6
$ws.Cells.Item(122, 8).Value = 2331.4363
$ws.Cells.Item(122, 9).Value = 2081.9487
$ws.Cells.Item(122, 10).Value = 2939.5625
$ws.Cells.Item(122, 11).Value = 6245.8461
$ws.Cells.Item(122, 12).Value = 8818.6875
$ws.Cells.Item(122, 13).Value = -3795.8461
$ws.Cells.Item(122, 14).Value = -13718.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1928.7441
$ws.Cells.Item(3, 9).Value = 650.30304
$ws.Cells.Item(3, 10).Value = 6147.6
$ws.Cells.Item(3, 11).Value = 650.30304
$ws.Cells.Item(3, 12).Value = 6147.6
$ws.Cells.Item(3, 13).Value = -536.30304
$ws.Cells.Item(3, 14).Value = -6375.6
$ws.Cells.Item(64, 8).Value = 1540.5
$ws.Cells.Item(64, 9).Value = 705
$ws.Cells.Item(64, 10).Value = 1819
$ws.Cells.Item(64, 11).Value = 705
$ws.Cells.Item(64, 12).Value = 1819
$ws.Cells.Item(64, 13).Value = -480
$ws.Cells.Item(64, 14).Value = -2269
$ws.Cells.Item(67, 8).Value = 1540.5
$ws.Cells.Item(67, 9).Value = 705
$ws.Cells.Item(67, 10).Value = 1819
$ws.Cells.Item(67, 11).Value = 705
$ws.Cells.Item(67, 12).Value = 1819
$ws.Cells.Item(67, 13).Value = 75
$ws.Cells.Item(67, 14).Value = -3379
$ws.Cells.Item(105, 8).Value = 1835.8096
$ws.Cells.Item(105, 9).Value = 1835.8096
$ws.Cells.Item(105, 11).Value = 1835.8096
$ws.Cells.Item(105, 13).Value = -88.80960000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 306.39285
$ws.Cells.Item(7, 9).Value = 80.375
$ws.Cells.Item(7, 11).Value = 80.375
$ws.Cells.Item(7, 13).Value = 32.625
$ws.Cells.Item(31, 8).Value = 17845.025
$ws.Cells.Item(31, 9).Value = 10040.77
$ws.Cells.Item(31, 10).Value = 21747.154
$ws.Cells.Item(31, 11).Value = 10040.77
$ws.Cells.Item(31, 12).Value = 21747.154
$ws.Cells.Item(31, 13).Value = -9745.77
$ws.Cells.Item(31, 14).Value = -22337.154
$ws.Cells.Item(34, 8).Value = 17845.025
$ws.Cells.Item(34, 9).Value = 10040.77
$ws.Cells.Item(34, 10).Value = 21747.154
$ws.Cells.Item(34, 11).Value = 10040.77
$ws.Cells.Item(34, 12).Value = 21747.154
$ws.Cells.Item(34, 13).Value = -9838.77
$ws.Cells.Item(34, 14).Value = -22151.154
$ws.Cells.Item(62, 8).Value = 3432.6843
$ws.Cells.Item(62, 9).Value = 2963.818
$ws.Cells.Item(62, 11).Value = 2963.818
$ws.Cells.Item(62, 13).Value = -2339.818
$ws.Cells.Item(65, 8).Value = 3432.6843
$ws.Cells.Item(65, 9).Value = 2963.818
$ws.Cells.Item(65, 11).Value = 14819.09
$ws.Cells.Item(65, 13).Value = -11699.09
$ws.Cells.Item(68, 8).Value = 45000
$ws.Cells.Item(68, 10).Value = 45000
$ws.Cells.Item(68, 12).Value = 45000
$ws.Cells.Item(68, 14).Value = -46498
$ws.Cells.Item(71, 8).Value = 45000
$ws.Cells.Item(71, 10).Value = 45000
$ws.Cells.Item(71, 12).Value = 135000
$ws.Cells.Item(71, 14).Value = -142488
$ws.Cells.Item(122, 8).Value = 4742.7427
$ws.Cells.Item(122, 9).Value = 3342.72
$ws.Cells.Item(122, 10).Value = 8242.799999999999
$ws.Cells.Item(122, 11).Value = 10028.16
$ws.Cells.Item(122, 12).Value = 24728.4
$ws.Cells.Item(122, 13).Value = -7578.16
$ws.Cells.Item(122, 14).Value = -29628.4
$ws.Cells.Item(134, 9).Value = 1451.4445
$ws.Cells.Item(134, 10).Value = 83347250
$ws.Cells.Item(134, 11).Value = 4354.333500000001
$ws.Cells.Item(134, 12).Value = 250041750
$ws.Cells.Item(134, 13).Value = -1819.333500000001
$ws.Cells.Item(134, 14).Value = -250046820

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 37389656
$ws.Cells.Item(4, 9).Value = 86991800
$ws.Cells.Item(4, 11).Value = 260975400
$ws.Cells.Item(4, 13).Value = -260975288
$ws.Cells.Item(38, 8).Value = 268
$ws.Cells.Item(38, 10).Value = 430.33334
$ws.Cells.Item(38, 12).Value = 1291.00002
$ws.Cells.Item(38, 14).Value = -1985.00002
$ws.Cells.Item(132, 8).Value = 1213.4445
$ws.Cells.Item(132, 9).Value = 1041.8572
$ws.Cells.Item(132, 11).Value = 9376.7148
$ws.Cells.Item(132, 13).Value = -6846.7148
$ws.Cells.Item(137, 8).Value = 2045.1364
$ws.Cells.Item(137, 9).Value = 2032
$ws.Cells.Item(137, 10).Value = 2060.9
$ws.Cells.Item(137, 11).Value = 6096
$ws.Cells.Item(137, 12).Value = 6182.700000000001
$ws.Cells.Item(137, 13).Value = -996
$ws.Cells.Item(137, 14).Value = -16382.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 4777.1665
$ws.Cells.Item(97, 9).Value = 1370.75
$ws.Cells.Item(97, 10).Value = 11590
$ws.Cells.Item(97, 11).Value = 1370.75
$ws.Cells.Item(97, 12).Value = 11590
$ws.Cells.Item(97, 13).Value = -874.75
$ws.Cells.Item(97, 14).Value = -12582
$ws.Cells.Item(103, 8).Value = 96460.39999999999
$ws.Cells.Item(103, 10).Value = 96460.39999999999
$ws.Cells.Item(103, 12).Value = 96460.39999999999
$ws.Cells.Item(103, 14).Value = -98804.39999999999
$ws.Cells.Item(122, 8).Value = 5851.963
$ws.Cells.Item(122, 9).Value = 3869.875
$ws.Cells.Item(122, 11).Value = 11609.625
$ws.Cells.Item(122, 13).Value = -9159.625
$ws.Cells.Item(132, 8).Value = 6427.871
$ws.Cells.Item(132, 9).Value = 2107.6667
$ws.Cells.Item(132, 11).Value = 6323.000100000001
$ws.Cells.Item(132, 13).Value = -3793.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 13393.5625
$ws.Cells.Item(22, 9).Value = 12754.272
$ws.Cells.Item(22, 10).Value = 14800
$ws.Cells.Item(22, 11).Value = 12754.272
$ws.Cells.Item(22, 12).Value = 14800
$ws.Cells.Item(22, 13).Value = -12459.272
$ws.Cells.Item(22, 14).Value = -15390
$ws.Cells.Item(27, 8).Value = 13393.5625
$ws.Cells.Item(27, 9).Value = 12754.272
$ws.Cells.Item(27, 10).Value = 14800
$ws.Cells.Item(27, 11).Value = 12754.272
$ws.Cells.Item(27, 12).Value = 14800
$ws.Cells.Item(27, 13).Value = -12647.272
$ws.Cells.Item(27, 14).Value = -15014
$ws.Cells.Item(40, 8).Value = 11954.363
$ws.Cells.Item(40, 9).Value = 7642.7144
$ws.Cells.Item(40, 11).Value = 7642.7144
$ws.Cells.Item(40, 13).Value = -7506.7144
$ws.Cells.Item(55, 8).Value = 185.85715
$ws.Cells.Item(55, 10).Value = 216.2
$ws.Cells.Item(55, 12).Value = 216.2
$ws.Cells.Item(55, 14).Value = -562.2
$ws.Cells.Item(93, 8).Value = 11624.263
$ws.Cells.Item(93, 9).Value = 7589.727
$ws.Cells.Item(93, 10).Value = 17171.75
$ws.Cells.Item(93, 11).Value = 7589.727
$ws.Cells.Item(93, 12).Value = 17171.75
$ws.Cells.Item(93, 13).Value = -6341.727
$ws.Cells.Item(93, 14).Value = -19667.75
$ws.Cells.Item(120, 8).Value = 73000
$ws.Cells.Item(120, 10).Value = 73000
$ws.Cells.Item(120, 12).Value = 73000
$ws.Cells.Item(120, 14).Value = -82676
$ws.Cells.Item(122, 8).Value = 6643.8887
$ws.Cells.Item(122, 9).Value = 4739.5835
$ws.Cells.Item(122, 10).Value = 10452.5
$ws.Cells.Item(122, 11).Value = 14218.7505
$ws.Cells.Item(122, 12).Value = 31357.5
$ws.Cells.Item(122, 13).Value = -11768.7505
$ws.Cells.Item(122, 14).Value = -36257.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 9981.75
$ws.Cells.Item(43, 9).Value = 9981.75
$ws.Cells.Item(43, 11).Value = 9981.75
$ws.Cells.Item(43, 13).Value = -9832.75
$ws.Cells.Item(107, 8).Value = 3675.8948
$ws.Cells.Item(107, 9).Value = 989.3333
$ws.Cells.Item(107, 10).Value = 13750.5
$ws.Cells.Item(107, 11).Value = 2967.9999
$ws.Cells.Item(107, 12).Value = 41251.5
$ws.Cells.Item(107, 13).Value = -1047.9999
$ws.Cells.Item(107, 14).Value = -45091.5
$ws.Cells.Item(122, 8).Value = 3184.3
$ws.Cells.Item(122, 10).Value = 8876.385
$ws.Cells.Item(122, 12).Value = 26629.155
$ws.Cells.Item(122, 14).Value = -31529.155
$ws.Cells.Item(126, 8).Value = 11878.116
$ws.Cells.Item(126, 9).Value = 13403.212
$ws.Cells.Item(126, 10).Value = 6845.3
$ws.Cells.Item(126, 11).Value = 40209.636
$ws.Cells.Item(126, 12).Value = 20535.9
$ws.Cells.Item(126, 13).Value = -37739.636
$ws.Cells.Item(126, 14).Value = -25475.9
$ws.Cells.Item(132, 8).Value = 5018.922
$ws.Cells.Item(132, 9).Value = 1866.8422
$ws.Cells.Item(132, 10).Value = 9625.808000000001
$ws.Cells.Item(132, 11).Value = 5600.5266
$ws.Cells.Item(132, 12).Value = 28877.424
$ws.Cells.Item(132, 13).Value = -3070.5266
$ws.Cells.Item(132, 14).Value = -33937.424
$ws.Cells.Item(139, 8).Value = 85999.5
$ws.Cells.Item(139, 10).Value = 85999.5
$ws.Cells.Item(139, 12).Value = 85999.5
$ws.Cells.Item(139, 14).Value = -96279.5
